# Apply the "term 2.0.0" update to KLMatterOfInterestValues.xlsx
#  - bump Version / Date / Contact / fix Description typo on the Metadata sheet
#  - insert a new Concept row (new UUID) at the top of the "Include from FSIII"
#    concept block, pushing the existing "B6" concept (and its following blank
#    row / System URI row) down by one row

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Metadata sheet
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B3").Value  = "2.0.0"                                                       # Version
$meta.Range("B8").Value  = "2024-06-04T14:59:10+02:00"                                   # Date
$meta.Range("B10").Value = "Kommunernes Landsforening (http://kl.dk)"                    # Contact
$meta.Range("B11").Value = "Matter of interest values to support when no observations have been made"  # Description (typo fix)

# ---------------------------------------------------------------------------
# 2. "Include from FSIII" sheet - add a new concept (d7ff926a-...) above "B6"
# ---------------------------------------------------------------------------
$fsiii = $wb.Worksheets.Item("Include from FSIII")

# Make room for the extra row: copy the formatting of the current last row
# (row 4, "System URI" / uri) down into the new row 5.
$fsiii.Range("A4:B4").Copy()
$fsiii.Range("A5:B5").PasteSpecial(-4122)

# Shift the existing three data rows down by one (bottom-up so we don't
# clobber a row before it has been read).
$fsiii.Range("A5").Value = $fsiii.Range("A4").Text
$fsiii.Range("B5").Value = $fsiii.Range("B4").Text

$fsiii.Range("A4").Value = $fsiii.Range("A3").Text
$fsiii.Range("B4").Value = $fsiii.Range("B3").Text

$fsiii.Range("A3").Value = $fsiii.Range("A2").Text
$fsiii.Range("B3").Value = $fsiii.Range("B2").Text

# New concept row
$fsiii.Range("A2").Value = "d7ff926a-4955-478f-b300-0b0ec0785013"
$fsiii.Range("B2").Value = ""
